# Saldo_guide.xlsx update: refresh extract date (2024-07-22 -> 2024-07-23),
# remove the client balance record for "CRISTINA FREIRE LIMA CARVALHO"
# (row 270, account 86633), and bump the sheet/window metadata to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the CRISTINA FREIRE LIMA CARVALHO row. It is row 270 in the
#    original sheet (B270=86633, E270=693.45, F270=29148219134). Deleting it
#    shifts every following row up by one and drops the now-unused shared
#    string from the workbook automatically.
$ws.Rows(270).Delete()

# 2) Every remaining data row's "Dt. Referencia" (column G) moves forward a
#    day, from 45495 (2024-07-22) to 45496 (2024-07-23). The sheet now has
#    274 data rows (2..274) after the deletion above.
for ($r = 2; $r -le 274; $r++) {
    $ws.Cells.Item($r, 7).Value = 45496
}

# 3) Rename the sheet to reflect the new extraction timestamp.
$ws.Name = "IClientBalance-20240723-094719-"

# 4) Reset the active selection back to A1 (was I14 in the source file).
[void]$ws.Range("A1").Select()
